$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (index 1): refresh the "想去人数" (want-to-go) counter column F
# for a batch of already-listed events. No rows are added or removed here.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 12752
$ws1.Range("F6").Value  = 53
$ws1.Range("F7").Value  = 35
$ws1.Range("F9").Value  = 6
$ws1.Range("F10").Value = 12650
$ws1.Range("F11").Value = 263
$ws1.Range("F12").Value = 15
$ws1.Range("F13").Value = 5902
$ws1.Range("F14").Value = 6537
$ws1.Range("F15").Value = 176
$ws1.Range("F16").Value = 83
$ws1.Range("F19").Value = 972
$ws1.Range("F23").Value = 180
$ws1.Range("F24").Value = 8
$ws1.Range("F25").Value = 82

# ---------------------------------------------------------------------------
# Sheet "演出" (index 2): append a brand-new event row (row 3) — the sheet
# previously only had the header row plus one data row.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Match the existing index-column formatting (bold, bordered, centred) before
# filling it in, so the new row looks like every other data row.
$ws2.Range("A2").Copy()
$ws2.Range("A3").PasteSpecial(-4122)
$ws2.Range("A3").Value = 2

# Column B holds plain "YYYY-MM-DD" text in this sheet (not real dates), so
# force text formatting before assigning or Excel will silently reinterpret
# it as a date serial number; the format is reset back to normal afterwards.
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "2024-10-15"
$ws2.Range("C1").Copy()
$ws2.Range("B3").PasteSpecial(-4122)

$ws2.Range("C3").Value = "苏州·Luca Stricagnoli 2024《进化时间》指弹吉他音乐会"
$ws2.Range("D3").Value = "滨河路999号红唐购物中心西区一层 山丘livehouse红唐店"
$ws2.Range("E3").Value = "2024.10.15 20:00-10.15 21:30"
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 220
$ws2.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=91359"
$ws2.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202408/9ncR7Aaj1724744018746.jpeg"

# ---------------------------------------------------------------------------
# Sheet "全部类型" (index 4): same counter refresh as sheet 1 (independent
# scrape pass, so a couple of values differ slightly), plus the same new
# "Luca Stricagnoli" event inserted in date order at row 23 — which pushes
# the previously-last five rows (old 23..27) down to 24..28.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F4").Value  = 12752
$ws4.Range("F7").Value  = 53
$ws4.Range("F8").Value  = 35
$ws4.Range("F10").Value = 6
$ws4.Range("F11").Value = 12650
$ws4.Range("F12").Value = 263
$ws4.Range("F13").Value = 15
$ws4.Range("F14").Value = 5904
$ws4.Range("F15").Value = 6538
$ws4.Range("F16").Value = 176
$ws4.Range("F17").Value = 83
$ws4.Range("F20").Value = 972

# Column B again needs the text guard for the whole block we are about to
# rewrite (old rows 23..27 shifting to 24..28, plus the new row 23).
$ws4.Range("B23:B28").NumberFormat = "@"

# Shift the tail of the table down by one row (bottom-up so nothing is
# clobbered before it is read), carrying the index-column formatting along.
$ws4.Range("A27").Copy()
$ws4.Range("A28").PasteSpecial(-4122)
$ws4.Range("A28").Value = 27
$ws4.Range("B28").Value = $ws4.Range("B27").Value2
$ws4.Range("C28").Value = $ws4.Range("C27").Value2
$ws4.Range("D28").Value = $ws4.Range("D27").Value2
$ws4.Range("E28").Value = $ws4.Range("E27").Value2
$ws4.Range("F28").Value = $ws4.Range("F27").Value2
$ws4.Range("G28").Value = $ws4.Range("G27").Value2
$ws4.Range("H28").Value = $ws4.Range("H27").Value2
$ws4.Range("I28").Value = $ws4.Range("I27").Value2

$ws4.Range("A27").Value = 26
$ws4.Range("B27").Value = $ws4.Range("B26").Value2
$ws4.Range("C27").Value = $ws4.Range("C26").Value2
$ws4.Range("D27").Value = $ws4.Range("D26").Value2
$ws4.Range("E27").Value = $ws4.Range("E26").Value2
$ws4.Range("F27").Value = $ws4.Range("F26").Value2
$ws4.Range("G27").Value = $ws4.Range("G26").Value2
$ws4.Range("H27").Value = $ws4.Range("H26").Value2
$ws4.Range("I27").Value = $ws4.Range("I26").Value2

$ws4.Range("A26").Value = 25
$ws4.Range("B26").Value = $ws4.Range("B25").Value2
$ws4.Range("C26").Value = $ws4.Range("C25").Value2
$ws4.Range("D26").Value = $ws4.Range("D25").Value2
$ws4.Range("E26").Value = $ws4.Range("E25").Value2
$ws4.Range("F26").Value = $ws4.Range("F25").Value2
$ws4.Range("G26").Value = $ws4.Range("G25").Value2
$ws4.Range("H26").Value = $ws4.Range("H25").Value2
$ws4.Range("I26").Value = $ws4.Range("I25").Value2

$ws4.Range("A25").Value = 24
$ws4.Range("B25").Value = $ws4.Range("B24").Value2
$ws4.Range("C25").Value = $ws4.Range("C24").Value2
$ws4.Range("D25").Value = $ws4.Range("D24").Value2
$ws4.Range("E25").Value = $ws4.Range("E24").Value2
$ws4.Range("F25").Value = $ws4.Range("F24").Value2
$ws4.Range("G25").Value = $ws4.Range("G24").Value2
$ws4.Range("H25").Value = $ws4.Range("H24").Value2
$ws4.Range("I25").Value = $ws4.Range("I24").Value2

$ws4.Range("A24").Value = 23
$ws4.Range("B24").Value = $ws4.Range("B23").Value2
$ws4.Range("C24").Value = $ws4.Range("C23").Value2
$ws4.Range("D24").Value = $ws4.Range("D23").Value2
$ws4.Range("E24").Value = $ws4.Range("E23").Value2
$ws4.Range("F24").Value = $ws4.Range("F23").Value2
$ws4.Range("G24").Value = $ws4.Range("G23").Value2
$ws4.Range("H24").Value = $ws4.Range("H23").Value2
$ws4.Range("I24").Value = $ws4.Range("I23").Value2

# Row 23 itself now becomes the new Luca Stricagnoli event (index value 22,
# same as it was before the shift — only the event details change).
$ws4.Range("A23").Value = 22
$ws4.Range("B23").Value = "2024-10-15"
$ws4.Range("C23").Value = "苏州·Luca Stricagnoli 2024《进化时间》指弹吉他音乐会"
$ws4.Range("D23").Value = "滨河路999号红唐购物中心西区一层 山丘livehouse红唐店"
$ws4.Range("E23").Value = "2024.10.15 20:00-10.15 21:30"
$ws4.Range("F23").Value = 0
$ws4.Range("G23").Value = 220
$ws4.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=91359"
$ws4.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202408/9ncR7Aaj1724744018746.jpeg"

# Reset column B's number format for the touched block back to the sheet's
# normal (unformatted) style now that the text values are safely stored.
$ws4.Range("C1").Copy()
$ws4.Range("B23:B28").PasteSpecial(-4122)

Write-Output "done"
